$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "D3"   = -7.560999999999998
    "E3"   = 16.987
    "C12"  = -10.771
    "D14"  = -7.431
    "D26"  = -7.639999999999999
    "E30"  = 16.323
    "D31"  = -7.915000000000001
    "C32"  = -13.407
    "D35"  = -7.935
    "C36"  = -12.721
    "D37"  = -7.712999999999999
    "C38"  = -12.703
    "E44"  = 16.613
    "D45"  = -7.603
    "C46"  = -14.015
    "C54"  = -12.705
    "C55"  = -13.624
    "E58"  = 16.624
    "C67"  = -11.616
    "C69"  = -10.641
    "C72"  = -11.555
    "E84"  = 16.354
    "E89"  = 17.199
    "C91"  = -11.169
    "E91"  = 16.987
    "E92"  = 16.971
    "C99"  = -12.635
    "D100" = -7.959000000000001
    "D102" = -7.695
    "E102" = 16.575
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
